# ---------------------------------------------------------------------------
# Day5_Rule_Based_Models_II.pptx -- "Update slides and workshops."
#
# 1) Bump every baked-in date field (datetimeFigureOut) from 7/26/2018 to
#    7/27/2018 across the slide master, every slide layout, and the notes
#    master.
# 2) On the last slide (closing "Reminders" slide):
#      - bold the "Optional office hours " / "4:30-5:30pm above Wu" spans
#        while leaving "today from " un-bold
#      - add a new bullet "AI4ALL_NLP_student -> git pull" after the
#        "Have a good weekend and take a break!" line (with a blank bullet
#        line in between)
#      - shrink the autofit font scale so the extra line still fits
#        (70% -> 55%)
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholders: 7/26/2018 -> 7/27/2018
#    ppPlaceholderDate == 16 identifies the "Date Placeholder" shape on the
#    slide master, every custom layout, and the notes master.
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.HasTextFrame) {
            $isDate = $false
            try {
                if ($shape.PlaceholderFormat.Type -eq 16) {
                    $isDate = $true
                }
            } catch {
                $isDate = $false
            }
            if ($isDate) {
                $tr = $shape.TextFrame.TextRange
                if ($tr.Text -eq "7/26/2018") {
                    $tr.Text = "7/27/2018"
                }
            }
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout on the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Notes master
Update-DatePlaceholder $p.NotesMaster.Shapes

# ---------------------------------------------------------------------------
# 2) Last slide ("Reminders") content tweaks
# ---------------------------------------------------------------------------
$lastSlide = $p.Slides.Item($p.Slides.Count)
$shape = $lastSlide.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

# --- Bold split of the "Optional office hours ..." bullet -----------------
$officeParaIdx = -1
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    if ($tr.Paragraphs($i).Text -eq "Optional office hours today from 4:30-5:30pm above Wu") {
        $officeParaIdx = $i
    }
}
$officePara = $tr.Paragraphs($officeParaIdx)
$officePara.Characters(1, 22).Font.Bold = $true   # "Optional office hours "
$officePara.Characters(34, 20).Font.Bold = $true  # "4:30-5:30pm above Wu"

# --- Add the new "AI4ALL_NLP_student -> git pull" bullet -------------------
$breakParaIdx = -1
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    if ($tr.Paragraphs($i).Text -eq "Have a good weekend and take a break!") {
        $breakParaIdx = $i
    }
}
$breakPara = $tr.Paragraphs($breakParaIdx)
$null = $breakPara.InsertAfter("`r`rAI4ALL_NLP_student -> git pull")
